$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each tuple: Datum, Receipt Number, Konto, Beskrivning, Debet, Kredit
# A value of $null means the source cell is blank (stored as an empty inline string).
# A string starting with a single quote forces Excel to store it as text, even when
# the text looks like a plain number (e.g. the "2261137" receipt number).
$data = @(
    ,(44368, $null, 6570, "Pris banktjänster enligt faktura", 97.5, $null)
    ,(44368, $null, $null, "Pris banktjänster enligt faktura", 0, $null)
    ,(44368, $null, 1930, "Pris banktjänster enligt faktura", $null, 97.5)
    ,(44368, "Reko233", 3011, "Reko Swish +46703533270", $null, 230.36)
    ,(44368, "Reko233", 2611, "Reko Swish +46703533270", $null, 27.64)
    ,(44368, "Reko233", 1930, "Reko Swish +46703533270", 258, $null)
    ,(44368, "Reko234", 3011, "Reko Swish +46709696209", $null, 115.18)
    ,(44368, "Reko234", 2611, "Reko Swish +46709696209", $null, 13.82)
    ,(44368, "Reko234", 1930, "Reko Swish +46709696209", 129, $null)
    ,(44368, "Reko235", 3011, "Reko Swish +46737812422", $null, 141.07)
    ,(44368, "Reko235", 2611, "Reko Swish +46737812422", $null, 16.93)
    ,(44368, "Reko235", 1930, "Reko Swish +46737812422", 158, $null)
    ,(44368, "Reko236", 3011, "Reko Swish +46767736486", $null, 742.86)
    ,(44368, "Reko236", 2611, "Reko Swish +46767736486", $null, 89.14)
    ,(44368, "Reko236", 1930, "Reko Swish +46767736486", 832, $null)
    ,(44368, "Reko237", 3011, "Reko Swish +46704105733", $null, 230.36)
    ,(44368, "Reko237", 2611, "Reko Swish +46704105733", $null, 27.64)
    ,(44368, "Reko237", 1930, "Reko Swish +46704105733", 258, $null)
    ,(44369, "Reko238", 3011, "Reko Swish +46767137127", $null, 691.0700000000001)
    ,(44369, "Reko238", 2611, "Reko Swish +46767137127", $null, 82.93000000000001)
    ,(44369, "Reko238", 1930, "Reko Swish +46767137127", 774, $null)
    ,(44369, "Reko239", 3011, "Reko Swish +46702174945", $null, 345.54)
    ,(44369, "Reko239", 2611, "Reko Swish +46702174945", $null, 41.46)
    ,(44369, "Reko239", 1930, "Reko Swish +46702174945", 387, $null)
    ,(44369, "Reko240", 3011, "Reko Swish +46739968231", $null, 460.71)
    ,(44369, "Reko240", 2611, "Reko Swish +46739968231", $null, 55.29)
    ,(44369, "Reko240", 1930, "Reko Swish +46739968231", 516, $null)
    ,(44369, "Reko241", 3011, "Reko Swish +46702453366", $null, 230.36)
    ,(44369, "Reko241", 2611, "Reko Swish +46702453366", $null, 27.64)
    ,(44369, "Reko241", 1930, "Reko Swish +46702453366", 258, $null)
    ,(44369, "Reko242", 3011, "Reko Swish +46761170330", $null, 34.82)
    ,(44369, "Reko242", 2611, "Reko Swish +46761170330", $null, 4.18)
    ,(44369, "Reko242", 1930, "Reko Swish +46761170330", 39, $null)
    ,(44369, "Reko243", 3011, "Reko Swish +46761910051", $null, 115.18)
    ,(44369, "Reko243", 2611, "Reko Swish +46761910051", $null, 13.82)
    ,(44369, "Reko243", 1930, "Reko Swish +46761910051", 129, $null)
    ,(44369, "Reko244", 3011, "Reko Swish +46768674881", $null, 70.54000000000001)
    ,(44369, "Reko244", 2611, "Reko Swish +46768674881", $null, 8.460000000000001)
    ,(44369, "Reko244", 1930, "Reko Swish +46768674881", 79, $null)
    ,(44369, "Reko245", 3011, "Reko Swish +46734029350", $null, 25.89)
    ,(44369, "Reko245", 2611, "Reko Swish +46734029350", $null, 3.11)
    ,(44369, "Reko245", 1930, "Reko Swish +46734029350", 29, $null)
    ,(44369, "Reko246", 3011, "Reko Swish +46705832242", $null, 141.07)
    ,(44369, "Reko246", 2611, "Reko Swish +46705832242", $null, 16.93)
    ,(44369, "Reko246", 1930, "Reko Swish +46705832242", 158, $null)
    ,(44369, "Reko247", 3011, "Reko Swish +46736699903", $null, 397.32)
    ,(44369, "Reko247", 2611, "Reko Swish +46736699903", $null, 47.68)
    ,(44369, "Reko247", 1930, "Reko Swish +46736699903", 445, $null)
    ,(44369, "Reko248", 3011, "Reko Swish +46737080200", $null, 115.18)
    ,(44369, "Reko248", 2611, "Reko Swish +46737080200", $null, 13.82)
    ,(44369, "Reko248", 1930, "Reko Swish +46737080200", 129, $null)
    ,(44369, "Reko249", 4010, "Reko Swish +46768674881 Return", 8.93, $null)
    ,(44369, "Reko249", 2645, "Reko Swish +46768674881 Return", 1.07, $null)
    ,(44369, "Reko249", 1930, "Reko Swish +46768674881 Return", $null, 10)
    ,(44370, $null, 5010, "Jun hyra", 4166, $null)
    ,(44370, $null, $null, "Jun hyra", 0, $null)
    ,(44370, $null, 1930, "Jun hyra", $null, 4166)
    ,(44373, $null, 4010, "NGROCERIES K0135", 334.82, $null)
    ,(44373, $null, 2645, "NGROCERIES K0135", 40.18, $null)
    ,(44373, $null, 1930, "NGROCERIES K0135", $null, 375)
    ,(44373, "'2261137", 3011, "Order 2261137 Swish +46735011685", $null, 530.36)
    ,(44373, "'2261137", 2611, "Order 2261137 Swish +46735011685", $null, 63.64)
    ,(44373, "'2261137", 1930, "Order 2261137 Swish +46735011685", 594, $null)
)

$startRow = 1268
$templateRow = 1267   # existing fully-populated data row whose formatting we reuse

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    # Copy number formats / styles from the last existing data row onto the new row
    $ws.Range("A$templateRow`:F$templateRow").Copy()
    $ws.Range("A$r`:F$r").PasteSpecial(-4122)

    for ($c = 0; $c -lt 6; $c++) {
        $v = $row[$c]
        $cell = $ws.Cells.Item($r, $c + 1)
        if ($v -eq $null) {
            $cell.ClearContents()
        } else {
            $cell.Value = $v
        }
    }
}

$excel.CutCopyMode = 0
